# LOQ4023.xlsx — site rebuild, 2023-04-12
#
# The course-info table had several rows whose "value" columns (B/C) were
# shifted up by one relative to their "label" column (A): the teacher's
# name was sitting next to "Objetivos:"/"Metodo:" instead of next to
# "Docentes responsaveis:", "Programa resumido:"/"Programa:" were showing
# a date/"Semestral" instead of the actual syllabus text, and
# "Bibliografia:" showed the make-up-exam rule. This fixes the layout by
# adding the missing value row and correcting each mis-placed cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the row that holds the "Docentes responsaveis:" value (pushes the
# old rows 13-24 down to 14-25 — each one keeps its own row height).
$ws.Rows.Item(13).Insert()

# The inserted row's A cell isn't used (A12 "Docentes responsaveis:" already
# spans down visually) — drop it so row 13 only carries B/C like the source.
$ws.Cells.Item(13, 1).Clear()

# Give the new B13/C13 cells the same look as the rest of the B/C columns
# (non-bold, wrapped text / red wrapped text) before filling them in.
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item(2, 3).Copy()
$ws.Cells.Item(13, 3).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Cells.Item(13, 2).Value = "1285870 - Marcos Villela Barcza"
$ws.Cells.Item(13, 3).Value = "1285870 - Marcos Villela Barcza"

# Row 10 "Objetivos:" — proper Portuguese objectives text (was the
# teacher's name).
$objetivos = "Conferir aos alunos uma noção prática das aplicações, à escala industrial, de processos e produtos estudados em disciplinas de química orgânica e ainda uma visão global das matérias primas mais importantes na área da indústria química de base orgânica."
$ws.Cells.Item(10, 2).Value = $objetivos
$ws.Cells.Item(10, 3).Value = $objetivos

# Row 14 (old 13) "Programa resumido:" — short Portuguese syllabus (was
# "Semestral").
$programaResumido = "Petróleo, Gás Natural e Petroquímica; Química Fina; Processos Unitários Orgânicos: Nitração; Esterificação; Alquilação e Acilação; Hidrogenação; Sulfonação/Sulfatação; Oxidação."
$ws.Cells.Item(14, 2).Value = $programaResumido
$ws.Cells.Item(14, 3).Value = $programaResumido

# Row 16 (old 15) "Programa:" — full Portuguese syllabus (was a leftover
# activation date).
$programa = "Petróleo, Gás Natural e Petroquímica; 2- Química Fina: Características, Química Fina X Química de Base, Principais Segmentos (Defensivos Agrícolas, Fármacos, Catalisadores, Corantes e Pigmentos, Especialidades); 3- Processos Unitários Orgânicos: 3.1- Nitração; 3.2- Esterificação; 3.3- Alquilação e Acilação; 3.4- Hidrogenação; 3.5- Sulfonação/Sulfatação; 3.6- Oxidação."
$ws.Cells.Item(16, 2).Value = $programa
$ws.Cells.Item(16, 3).Value = $programa

# Row 19 (old 18) "Metodo:" — teaching method text (was the teacher's
# name again).
$metodo = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos."
$ws.Cells.Item(19, 2).Value = $metodo
$ws.Cells.Item(19, 3).Value = $metodo

# Row 20 (old 19) "Criterio:" — grading criteria (was the method text).
$criterio = "Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula."
$ws.Cells.Item(20, 2).Value = $criterio
$ws.Cells.Item(20, 3).Value = $criterio

# Row 21 (old 20) "Norma de recuperacao:" — make-up exam rule (was the
# grading-criteria text).
$normaRecuperacao = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
$ws.Cells.Item(21, 2).Value = $normaRecuperacao
$ws.Cells.Item(21, 3).Value = $normaRecuperacao

# Row 22 (old 21) "Bibliografia:" — the actual bibliography (was the
# make-up-exam rule).
$bibliografia = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;Química & Derivados, São Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Disponível em: http://www.quimica.com.br/pquimica/category/revista/BiodieselBR. Curitiba: BiodieselBR, v.1, n.1, out/dez.2007 -;Petróleo & Energia, São Paulo, SP: , v. 1, n. 1, ; Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/."
$ws.Cells.Item(22, 2).Value = $bibliografia
$ws.Cells.Item(22, 3).Value = $bibliografia

# Column layout: column A now gets its own (narrower) width instead of
# sharing a merged 1-2 column-width definition with column B.
$ws.Columns.Item(1).ColumnWidth = 30.71

Write-Host "LOQ4023 table corrected"
